$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.527.09'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '2.482.14'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.67'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  -1.35%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.508'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.111'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.20%  '
$ws.Range('D13').Value = '2.866.01'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('E14').Value = '  -1.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.08'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +10.50%  '
$ws.Range('D16').Value = '2.531.38'
$ws.Range('E16').Value = '  +3.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.764'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.48%  '
$ws.Range('D18').Value = '41.532.07'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('E19').Value = '  +2.42%  '
$ws.Range('D20').Value = '0.0₃0934'
$ws.Range('E20').Value = '  +2.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.36'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('E26').Value = '  +0.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.94'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.19'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.05'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.78%  '
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0754'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.09%  '
$ws.Range('E35').Value = '  -7.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.21%  '
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.104'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.58%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.93'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').Value = '1.977.69'
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.96'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.12'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +7.58%  '
$ws.Range('D48').Value = '2.722.58'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.58'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.05%  '
